$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated GDP per Capita values by year (stored as text, matching source sheet formatting)
$values = @{
    1950 = "902"
    1951 = "942"
    1952 = "958"
    1953 = "982"
    1954 = "1023"
    1955 = "1044"
    1956 = "1070"
    1957 = "1093"
    1958 = "1109"
    1959 = "1157"
    1960 = "1176"
    1961 = "1181"
    1962 = "1216"
    1963 = "1218"
    1964 = "1215"
    1965 = "1282"
    1966 = "1342"
    1967 = "1313"
    1968 = "1420"
    1969 = "1525"
    1970 = "1553"
    1971 = "1535"
    1972 = "1525"
    1973 = "1553"
    1974 = "1581"
    1975 = "1551"
    1976 = "1511"
    1977 = "1495"
    1978 = "1524"
    1979 = "1572"
    1980 = "1584"
    1981 = "1615"
    1982 = "1610"
    1983 = "1575"
    1984 = "1608"
    1985 = "1610"
    1986 = "1597"
    1987 = "1624"
    1988 = "1626"
    1989 = "1642"
    1990 = "1635"
    1991 = "1467.38879372737"
    1992 = "1354.96819689212"
    1993 = "1405.16280529918"
    1994 = "1500.20920326603"
    1995 = "1383.10812635823"
    1996 = "1060.56242298435"
    1997 = "885.298443424568"
    1998 = "882.729341181711"
    1999 = "814.123267008335"
    2000 = "828.793864301684"
    2001 = "923.045794079924"
    2002 = "1101.75037870797"
    2003 = "1156.51911986097"
    2004 = "1197.15292990557"
    2005 = "1223.70153268588"
    2006 = "1247.25275616654"
    2007 = "1320.75766231695"
    2008 = "1363.48755538923"
    2009 = "1378.35793886063"
    2010 = "1421.15205985533"
    2011 = "1479"
    2012 = "1665"
    2013 = "1965"
    2014 = "2007"
    2015 = "1563"
    2016 = "1619"
}

# Existing data rows 2..60 correspond to years 1950..2008 -> refresh column E (Data)
for ($row = 2; $row -le 60; $row++) {
    $year = 1948 + $row
    $ws.Cells.Item($row, 5).Value = "'" + $values[$year]
}

# New rows for years 2009..2016 (rows 61..68)
for ($row = 61; $row -le 68; $row++) {
    $year = 1948 + $row
    $ws.Cells.Item($row, 1).Value = 694
    $ws.Cells.Item($row, 2).Value = "Sierra Leone"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = "'" + $values[$year]
}
